$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.210757114293816
$ws.Range("C2").Value = 0.210757114293816
$ws.Range("D2").Value = 0.101347349640299
$ws.Range("E2").Value = 0.000658751049003168
$ws.Range("F2").Value = 0.7554
$ws.Range("B3").Value = 19.9493093126419
$ws.Range("C3").Value = 19.9493093126419
$ws.Range("D3").Value = 9.59307889921189
$ws.Range("E3").Value = 0.0623543764139359
$ws.Range("F3").Value = 0.0019
$ws.Range("B4").Value = 0.318790907760633
$ws.Range("C4").Value = 0.318790907760633
$ws.Range("D4").Value = 0.153297855207506
$ws.Range("E4").Value = 0.000996425888652198
$ws.Range("F4").Value = 0.6988
$ws.Range("B5").Value = 299.455531555822
$ws.Range("C5").Value = 2.07955230247099
$ws.Range("E5").Value = 0.935990446648409
$ws.Range("B6").Value = 319.934388890519
$ws.Range("B7").Value = 0.0149852447850418
$ws.Range("C7").Value = 0.0149852447850418
$ws.Range("D7").Value = 0.0184668245087217
$ws.Range("E7").Value = 0.000127657468861562
$ws.Range("F7").Value = 0.89
$ws.Range("B8").Value = 0.0982469912216245
$ws.Range("C8").Value = 0.0982469912216245
$ws.Range("D8").Value = 0.121073093661486
$ws.Range("E8").Value = 0.000836954110695343
$ws.Range("F8").Value = 0.7316
$ws.Range("B9").Value = 0.42166987486296
$ws.Range("C9").Value = 0.42166987486296
$ws.Range("D9").Value = 0.519638063402327
$ws.Range("E9").Value = 0.00359215412843367
$ws.Range("F9").Value = 0.4697
$ws.Range("B10").Value = 116.851451532821
$ws.Range("C10").Value = 0.811468413422371
$ws.Range("E10").Value = 0.995443234292009
$ws.Range("B11").Value = 117.386353643691
$ws.Range("B12").Value = 0.0243295345328576
$ws.Range("C12").Value = 0.0243295345328576
$ws.Range("D12").Value = 0.0710358119133855
$ws.Range("E12").Value = 0.000485902785068253
$ws.Range("F12").Value = 0.7896
$ws.Range("B13").Value = 0.0395669751608717
$ws.Range("C13").Value = 0.0395669751608717
$ws.Range("D13").Value = 0.115525112151792
$ws.Range("E13").Value = 0.000790220766510316
$ws.Range("F13").Value = 0.738
$ws.Range("B14").Value = 0.687358196769216
$ws.Range("C14").Value = 0.687358196769216
$ws.Range("D14").Value = 2.00690430459653
$ws.Range("E14").Value = 0.0137277292213953
$ws.Range("F14").Value = 0.1604
$ws.Range("B15").Value = 49.3195316329076
$ws.Range("C15").Value = 0.342496747450747
$ws.Range("E15").Value = 0.984996147227026
$ws.Range("B16").Value = 50.0707863393706
$ws.Range("F17").Value = 0.9891
$ws.Range("F18").Value = 0.2173
$ws.Range("F19").Value = 0.7643
$ws.Range("F22").Value = 0.0857
$ws.Range("F23").Value = 0.0429
$ws.Range("F24").Value = 0.1401
$ws.Range("B27").Value = 0.0245046621802622
$ws.Range("C27").Value = 0.0245046621802622
$ws.Range("D27").Value = 0.0503732173487164
$ws.Range("E27").Value = 0.00034885943435843
$ws.Range("F27").Value = 0.8262
$ws.Range("B28").Value = 0.0507987323022581
$ws.Range("C28").Value = 0.0507987323022581
$ws.Range("D28").Value = 0.104424846360952
$ws.Range("E28").Value = 0.000723193688071543
$ws.Range("F28").Value = 0.7438
$ws.Range("B29").Value = 0.116374093419246
$ws.Range("C29").Value = 0.116374093419246
$ws.Range("D29").Value = 0.239225395495934
$ws.Range("E29").Value = 0.00165675413541974
$ws.Range("F29").Value = 0.634
$ws.Range("B30").Value = 70.0505455017888
$ws.Range("C30").Value = 0.4864621215402
$ws.Range("E30").Value = 0.99727119274215
$ws.Range("B31").Value = 70.2422229896906
$ws.Range("B32").Value = 0.0219846656846166
$ws.Range("C32").Value = 0.0219846656846166
$ws.Range("D32").Value = 0.018768248912626
$ws.Range("E32").Value = 0.000119723027914035
$ws.Range("F32").Value = 0.8901
$ws.Range("B33").Value = 14.9161045557755
$ws.Range("C33").Value = 14.9161045557755
$ws.Range("D33").Value = 12.7338376268982
$ws.Range("E33").Value = 0.0812293999698778
$ws.Range("F33").Value = 0.001
$ws.Range("B34").Value = 0.0132310001462539
$ws.Range("C34").Value = 0.0132310001462539
$ws.Range("D34").Value = 0.0112952686054099
$ws.Range("E34").Value = 0.0000720527399672477
$ws.Range("F34").Value = 0.9153
$ws.Range("B35").Value = 168.678062259451
$ws.Range("C35").Value = 1.1713754323573
$ws.Range("E35").Value = 0.918578824262241
$ws.Range("B36").Value = 183.629382481057
$ws.Range("F37").Value = 0.6821
$ws.Range("F38").Value = 0.0495
$ws.Range("F39").Value = 0.7337
